$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.22599999999999
$ws.Range("C14").Value = -12.3539
$ws.Range("C16").Value = -11.86580000000001
$ws.Range("C21").Value = -13.16310000000002
$ws.Range("C23").Value = -12.03090000000001
$ws.Range("C25").Value = -11.07519999999999
